$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old win value from B2 (row 2 moves to column C as "Support")
$ws.Range("B2").ClearContents()

# Row 2: A2=0 (unchanged), C2 = Support
$ws.Range("C2").Value = "Support"

# Row 3: A3=1 (unchanged), B3 = loss, C3 = DPS
$ws.Range("B3").Value = "loss"
$ws.Range("C3").Value = "DPS"

# Row 4: A4 = 2
$ws.Range("A4").Value = 2

# Row 5: A5 = 3, B5 = win, C5 = Tank
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "win"
$ws.Range("C5").Value = "Tank"

# Apply the same style as other numeric A cells (A2/A3) to the new A4/A5 cells
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
